# Today_words.xlsx — "Add another aprox 50 words on 20 JAN 2012"
# The existing three study words are cleared out (content removed, the
# blue "answer" style is kept on B2:B3) and the sheet is re-laid out as a
# two-column word list: column B keeps its old formatting, a new column C
# is added (Arial 10, black) with ~24 spacer/answer rows sprinkled down to
# row 60, ready for new words to be typed in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Remove the three existing shared-string words, keep the blue style
#    that was already applied to B2:B3.
$ws.Range("B1:B3").ClearContents()

# 2) New column for the second word list: same row height, own width
#    (~26.1 chars, matching the width Excel wrote for the new column).
$ws.Columns.Item(3).ColumnWidth = 25.25

# 3) Build the "answer" format (Arial 10, solid black) once on a scratch
#    cell, then fan it out (format only, no value) to every row that will
#    hold a word in column C.
$tmpl = $ws.Range("E1")
$tmpl.Font.Color = 0
$tmpl.Font.Size = 10
$tmpl.Font.Name = "Arial"
$tmpl.Copy()

$cRows = @(1,16,19,27,30,31,34,35,36,38,39,40,43,44,46,47,48,52,55,56,57,58,59,60)
$xlPasteFormats = -4122
foreach ($r in $cRows) {
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)
}
$tmpl.Clear()
$excel.CutCopyMode = $false

# 4) Printable, portrait page for the new layout.
$ws.PageSetup.Orientation = 1

# 5) Match the selection left behind by the edit.
$ws.Range("A1:C61").Select() | Out-Null
